$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")
$ws.Range("Z200").Formula = "='[Limits-V5-Input.xlsx]main'!D24"
$ws.Range("Z201").Formula = "='[Limits-V5-Input.xlsx]input-output'!D24"
$ws.Range("Z202").Formula = "='[Limits-V5-Input.xlsx]Some Calculations'!D24"
$ws.Range("Z203").Formula = "='[Limits-V5-Input.xlsx]Report'!D24"
$ws.Range("Z200:Z203").ClearContents()
